# Scheduled runner update: refresh market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N) across
# several sheets of the Asura_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2015.4166
$ws.Range("I62").Value = 2229
$ws.Range("J62").Value = 947.5
$ws.Range("K62").Value = 2229
$ws.Range("L62").Value = 947.5
$ws.Range("M62").Value = -1605
$ws.Range("N62").Value = -2195.5
$ws.Range("H65").Value = 2015.4166
$ws.Range("I65").Value = 2229
$ws.Range("J65").Value = 947.5
$ws.Range("K65").Value = 11145
$ws.Range("L65").Value = 4737.5
$ws.Range("M65").Value = -8025
$ws.Range("N65").Value = -10977.5
$ws.Range("H69").Value = 4430.4287
$ws.Range("I69").Value = 3671
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 11013
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -10139
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 4430.4287
$ws.Range("I72").Value = 3671
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 33039
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -28671
$ws.Range("N72").Value = -53736
$ws.Range("H96").Value = 987.3333
$ws.Range("I96").Value = 792.3333
$ws.Range("J96").Value = 1572.3334
$ws.Range("K96").Value = 2376.9999
$ws.Range("L96").Value = 4717.0002
$ws.Range("M96").Value = -1003.9999
$ws.Range("N96").Value = -7463.0002
$ws.Range("H100").Value = 3677.3076
$ws.Range("I100").Value = 3400.7144
$ws.Range("K100").Value = 3400.7144
$ws.Range("M100").Value = -2859.7144
$ws.Range("H132").Value = 1449.4286
$ws.Range("I132").Value = 1360.8276
$ws.Range("J132").Value = 2477.2
$ws.Range("K132").Value = 4082.4828
$ws.Range("L132").Value = 7431.599999999999
$ws.Range("M132").Value = -1552.4828
$ws.Range("N132").Value = -12491.6
$ws.Range("H138").Value = 2441.925
$ws.Range("I138").Value = 1397.6097
$ws.Range("J138").Value = 3539.795
$ws.Range("K138").Value = 4192.8291
$ws.Range("L138").Value = 10619.385
$ws.Range("M138").Value = 947.1709000000001
$ws.Range("N138").Value = -20899.385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26558.256
$ws.Range("I2").Value = 653.85
$ws.Range("K2").Value = 653.85
$ws.Range("M2").Value = -540.85
$ws.Range("H43").Value = 500008700
$ws.Range("J43").Value = 500008700
$ws.Range("L43").Value = 500008700
$ws.Range("N43").Value = -500009326
$ws.Range("H63").Value = 4442.857
$ws.Range("I63").Value = 3880
$ws.Range("J63").Value = 5850
$ws.Range("K63").Value = 3880
$ws.Range("L63").Value = 5850
$ws.Range("M63").Value = -3194
$ws.Range("N63").Value = -7222
$ws.Range("H66").Value = 4442.857
$ws.Range("I66").Value = 3880
$ws.Range("J66").Value = 5850
$ws.Range("K66").Value = 19400
$ws.Range("L66").Value = 29250
$ws.Range("M66").Value = -15968
$ws.Range("N66").Value = -36114
$ws.Range("H74").Value = 985.375
$ws.Range("I74").Value = 1002.6667
$ws.Range("J74").Value = 975
$ws.Range("K74").Value = 1002.6667
$ws.Range("L74").Value = 975
$ws.Range("M74").Value = -128.6667
$ws.Range("N74").Value = -2723
$ws.Range("H77").Value = 985.375
$ws.Range("I77").Value = 1002.6667
$ws.Range("J77").Value = 975
$ws.Range("K77").Value = 5013.3335
$ws.Range("L77").Value = 4875
$ws.Range("M77").Value = -645.3334999999997
$ws.Range("N77").Value = -13611
$ws.Range("H116").Value = 26558.256
$ws.Range("I116").Value = 653.85
$ws.Range("K116").Value = 653.85
$ws.Range("M116").Value = 1640.15
$ws.Range("H122").Value = 2783.3235
$ws.Range("I122").Value = 2541.923
$ws.Range("J122").Value = 3567.875
$ws.Range("K122").Value = 7625.768999999999
$ws.Range("L122").Value = 10703.625
$ws.Range("M122").Value = -5175.768999999999
$ws.Range("N122").Value = -15603.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26558.256
$ws.Range("I3").Value = 653.85
$ws.Range("K3").Value = 653.85
$ws.Range("M3").Value = -539.85
$ws.Range("H105").Value = 2642.9167
$ws.Range("I105").Value = 2446.818
$ws.Range("J105").Value = 4800
$ws.Range("K105").Value = 2446.818
$ws.Range("L105").Value = 4800
$ws.Range("M105").Value = -699.8180000000002
$ws.Range("N105").Value = -8294
$ws.Range("H134").Value = 2978.5356
$ws.Range("I134").Value = 2662.6875
$ws.Range("J134").Value = 3399.6667
$ws.Range("K134").Value = 7988.0625
$ws.Range("L134").Value = 10199.0001
$ws.Range("M134").Value = -5453.0625
$ws.Range("N134").Value = -15269.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2746.2
$ws.Range("I31").Value = 1459.7333
$ws.Range("J31").Value = 6605.6
$ws.Range("K31").Value = 1459.7333
$ws.Range("L31").Value = 6605.6
$ws.Range("M31").Value = -1164.7333
$ws.Range("N31").Value = -7195.6
$ws.Range("H34").Value = 2746.2
$ws.Range("I34").Value = 1459.7333
$ws.Range("J34").Value = 6605.6
$ws.Range("K34").Value = 1459.7333
$ws.Range("L34").Value = 6605.6
$ws.Range("M34").Value = -1257.7333
$ws.Range("N34").Value = -7009.6
$ws.Range("H105").Value = 2222.5
$ws.Range("I105").Value = 2222.5
$ws.Range("K105").Value = 2222.5
$ws.Range("M105").Value = -475.5
$ws.Range("H107").Value = 404.42856
$ws.Range("I107").Value = 371.41666
$ws.Range("J107").Value = 448.44446
$ws.Range("K107").Value = 371.41666
$ws.Range("L107").Value = 448.44446
$ws.Range("M107").Value = 1548.58334
$ws.Range("N107").Value = -4288.44446
$ws.Range("H134").Value = 1714.9062
$ws.Range("I134").Value = 1224.4783
$ws.Range("J134").Value = 2968.2222
$ws.Range("K134").Value = 3673.4349
$ws.Range("L134").Value = 8904.6666
$ws.Range("M134").Value = -1138.4349
$ws.Range("N134").Value = -13974.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 368.72726
$ws.Range("I107").Value = 406.55
$ws.Range("J107").Value = 310.53845
$ws.Range("K107").Value = 1219.65
$ws.Range("L107").Value = 931.61535
$ws.Range("M107").Value = 700.3499999999999
$ws.Range("N107").Value = -4771.61535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 96117.64999999999
$ws.Range("I97").Value = 69553.336
$ws.Range("J97").Value = 126002.5
$ws.Range("K97").Value = 69553.336
$ws.Range("L97").Value = 126002.5
$ws.Range("M97").Value = -69057.336
$ws.Range("N97").Value = -126994.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1324.475
$ws.Range("I93").Value = 698.03845
$ws.Range("J93").Value = 2487.8572
$ws.Range("K93").Value = 698.03845
$ws.Range("L93").Value = 2487.8572
$ws.Range("M93").Value = 549.96155
$ws.Range("N93").Value = -4983.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12500
$ws.Range("J45").Value = 12500
$ws.Range("L45").Value = 12500
$ws.Range("N45").Value = -13482
$ws.Range("H74").Value = 13333.333
$ws.Range("J74").Value = 13333.333
$ws.Range("L74").Value = 13333.333
$ws.Range("N74").Value = -15205.333
$ws.Range("H77").Value = 13333.333
$ws.Range("J77").Value = 13333.333
$ws.Range("L77").Value = 39999.999
$ws.Range("N77").Value = -49359.999
